# Insert a new weekly data row for "Puerro" (Vega Modelo de Temuco) above the
# existing row 280, shifting all subsequent rows (280-314 -> 281-315) down by
# one. Rows.Insert() takes care of re-indexing every row below, growing the
# sheet's used range / dimension from R314 to R315 automatically.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(280).Insert()

$ws.Range("A280").Value = 10
$ws.Range("B280").Value = "Vega Modelo de Temuco"
$ws.Range("C280").Value = "La Araucanía"
$ws.Range("D280").Value = 45124
$ws.Range("E280").Value = 9
$ws.Range("F280").Value = 100112005
$ws.Range("G280").Value = "Puerro"
$ws.Range("H280").Value = "Azul de Maquehue"
$ws.Range("I280").Value = "Primera"
$ws.Range("J280").Value = 50
$ws.Range("K280").Value = 8000
$ws.Range("L280").Value = 8000
$ws.Range("M280").Value = 8000
$ws.Range("N280").Value = "`$/docena de paquetes"
$ws.Range("O280").Value = "Provincia de Cautín"
$ws.Range("P280").Value = 667
$ws.Range("Q280").Value = 12
$ws.Range("R280").Value = "Hortaliza"
